# fix: contact transmission per cm2
# The "hand_surface" input row (Worker hand surface area, sourced from
# Google, 450 cm2) is removed from the inputs table. All rows below it
# shift up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 holds the hand_surface variable - delete the entire row, which
# shifts every subsequent row up by one and updates the shared strings /
# dimension automatically.
$ws.Rows.Item(11).Delete() | Out-Null

# Restore the cursor/selection position recorded in the saved workbook.
$ws.Range("D27").Select() | Out-Null
